$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.2269230769230769
$ws.Range("C2").Value = 0.4884615384615384
$ws.Range("J2").Value = 0.03076923076923077
$ws.Range("P2").Value = 0.1692307692307692
$ws.Range("S2").Value = 0.08461538461538462
$ws.Range("C3").Value = 0.03007518796992481
$ws.Range("J3").Value = 0.05263157894736842
$ws.Range("P3").Value = 0.7218045112781954
$ws.Range("S3").Value = 0.1954887218045113
$ws.Range("J4").Value = 0.03333333333333333
$ws.Range("P4").Value = 0.7333333333333333
$ws.Range("S4").Value = 0.2333333333333333
$ws.Range("B6").Value = 0.0423728813559322
$ws.Range("D6").Value = 0.00423728813559322
$ws.Range("F6").Value = 0.05084745762711865
$ws.Range("J6").Value = 0.2754237288135593
$ws.Range("O6").Value = 0.008474576271186441
$ws.Range("Q6").Value = 0.1694915254237288
$ws.Range("R6").Value = 0.08898305084745763
$ws.Range("S6").Value = 0.3601694915254237
$ws.Range("B7").Value = 0.1123595505617977
$ws.Range("D7").Value = 0.005617977528089887
$ws.Range("F7").Value = 0.05056179775280899
$ws.Range("J7").Value = 0.09550561797752809
$ws.Range("O7").Value = 0.02808988764044944
$ws.Range("Q7").Value = 0.1910112359550562
$ws.Range("R7").Value = 0.06741573033707865
$ws.Range("S7").Value = 0.449438202247191
$ws.Range("B8").Value = 0.09042553191489362
$ws.Range("D8").Value = 0.02393617021276596
$ws.Range("F8").Value = 0.07712765957446809
$ws.Range("J8").Value = 0.07712765957446809
$ws.Range("O8").Value = 0.01329787234042553
$ws.Range("Q8").Value = 0.1595744680851064
$ws.Range("R8").Value = 0.09308510638297872
$ws.Range("S8").Value = 0.4654255319148936
$ws.Range("B9").Value = 0.05172413793103448
$ws.Range("D9").Value = 0.01149425287356322
$ws.Range("F9").Value = 0.05747126436781609
$ws.Range("J9").Value = 0.09770114942528736
$ws.Range("O9").Value = 0.02298850574712644
$ws.Range("Q9").Value = 0.2241379310344828
$ws.Range("R9").Value = 0.08620689655172414
$ws.Range("S9").Value = 0.4482758620689655
$ws.Range("B10").Value = 0.1017087062652563
$ws.Range("D10").Value = 0.01301871440195281
$ws.Range("F10").Value = 0.07973962571196094
$ws.Range("J10").Value = 0.08624898291293735
$ws.Range("O10").Value = 0.01057770545158665
$ws.Range("Q10").Value = 0.2050447518307567
$ws.Range("R10").Value = 0.0870626525630594
$ws.Range("S10").Value = 0.4165988608624898
$ws.Range("G11").Value = 0.1262135922330097
$ws.Range("J11").Value = 0.116504854368932
$ws.Range("K11").Value = 0.1974110032362459
$ws.Range("L11").Value = 0.5339805825242718
$ws.Range("S11").Value = 0.02588996763754045
$ws.Range("G12").Value = 0.6964285714285714
$ws.Range("J12").Value = 0.244047619047619
$ws.Range("K12").Value = 0.02976190476190476
$ws.Range("L12").Value = 0.005952380952380952
$ws.Range("S12").Value = 0.02380952380952381
$ws.Range("G13").Value = 0.6486486486486487
$ws.Range("J13").Value = 0.2702702702702703
$ws.Range("S13").Value = 0.08108108108108109
$ws.Range("F15").Value = 0.01834862385321101
$ws.Range("H15").Value = 0.1192660550458716
$ws.Range("I15").Value = 0.1055045871559633
$ws.Range("J15").Value = 0.4082568807339449
$ws.Range("K15").Value = 0.07798165137614679
$ws.Range("M15").Value = 0.01376146788990826
$ws.Range("O15").Value = 0.06880733944954129
$ws.Range("S15").Value = 0.1880733944954129
$ws.Range("F16").Value = 0.03125
$ws.Range("H16").Value = 0.20625
$ws.Range("I16").Value = 0.06875000000000001
$ws.Range("J16").Value = 0.4
$ws.Range("K16").Value = 0.10625
$ws.Range("M16").Value = 0.00625
$ws.Range("N16").Value = 0.00625
$ws.Range("O16").Value = 0.05
$ws.Range("S16").Value = 0.125
$ws.Range("F17").Value = 0.02364066193853428
$ws.Range("H17").Value = 0.1773049645390071
$ws.Range("I17").Value = 0.06382978723404255
$ws.Range("J17").Value = 0.4231678486997636
$ws.Range("K17").Value = 0.1087470449172577
$ws.Range("M17").Value = 0.01182033096926714
$ws.Range("N17").Value = 0.002364066193853428
$ws.Range("O17").Value = 0.06382978723404255
$ws.Range("S17").Value = 0.1252955082742317
$ws.Range("F18").Value = 0.03191489361702127
$ws.Range("H18").Value = 0.101063829787234
$ws.Range("I18").Value = 0.0797872340425532
$ws.Range("J18").Value = 0.4787234042553192
$ws.Range("K18").Value = 0.09042553191489362
$ws.Range("M18").Value = 0.03723404255319149
$ws.Range("O18").Value = 0.06382978723404255
$ws.Range("S18").Value = 0.1170212765957447
$ws.Range("F19").Value = 0.02285263987391647
$ws.Range("H19").Value = 0.1773049645390071
$ws.Range("I19").Value = 0.07801418439716312
$ws.Range("J19").Value = 0.3877068557919622
$ws.Range("K19").Value = 0.1150512214342002
$ws.Range("M19").Value = 0.01654846335697399
$ws.Range("O19").Value = 0.0764381402679275
$ws.Range("S19").Value = 0.1260835303388495
